# "test(results.xlsx): add user studies 2 results"
#
# Fill in the "Nina User2" answers column (B, plus a couple of spillover
# cells in C/D on row 13) and make that sheet the active one with its
# used range selected - mirroring what the author did by hand after
# transcribing the second round of user-study answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nina User2")

$ws.Range("B2").Value  = "Cat"
$ws.Range("B3").Value  = "Obama"
$ws.Range("B4").Value  = "Eiffel Tower"
$ws.Range("B5").Value  = "Lady Gaga"
$ws.Range("B6").Value  = "Prince of England, Marriage"
$ws.Range("B7").Value  = "Two Grizzly bears"
$ws.Range("B8").Value  = "Gumball machine"
$ws.Range("B9").Value  = "The painting, with the pitchfork"
$ws.Range("B10").Value = "Yoga, outside in the mountains"
$ws.Range("B11").Value = "Some white dude cooking"
$ws.Range("B12").Value = "Starship"
$ws.Range("B13").Value = "Dragon, fire"
$ws.Range("C13").Value = "C5"
$ws.Range("D13").Value = "Eggrolls"
$ws.Range("B14").Value = "Horse racing"
$ws.Range("B16").Value = "Clocks, doorway"
$ws.Range("B17").Value = "Family doing kissy faces"
$ws.Range("B18").Value = "Graduation"
$ws.Range("B19").Value = "Nina and Molly sleeping"
$ws.Range("B20").Value = "Robert and Happy"
$ws.Range("B21").Value = "Happy and Molly"
$ws.Range("B15").Value = "Person in the bottom right, storefronts"

# Portrait page setup, matching the other sheets in this workbook.
$ws.PageSetup.Orientation = 1

# Make "Nina User2" the active sheet/tab, with its newly-filled range
# selected (A1:D21), same as sheet1 was tabSelected before the edit.
$ws.Activate()
$ws.Range("A1:D21").Select()
